$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.665.77'
$ws.Range('E2').Value = '  +4.53%  '
$ws.Range('D3').Value = '4.036.95'
$ws.Range('E3').Value = '  +4.54%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '534.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.08'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +7.78%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.693'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +13.86%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.758'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.37%  '
$ws.Range('E10').Value = '  +5.07%  '
$ws.Range('E11').Value = '  +3.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '48.54'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +16.91%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.91'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.48%  '
$ws.Range('D14').Value = '4.692.71'
$ws.Range('E14').Value = '  +5.05%  '
$ws.Range('D15').Value = '4.049.88'
$ws.Range('E15').Value = '  +5.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.36'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.44%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '20.92'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.56%  '
$ws.Range('E18').Value = '  +1.87%  '
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').Value = '71.731.63'
$ws.Range('E20').Value = '  +4.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '436.13'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.68'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '99.13'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +14.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.72'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +5.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.24'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.98%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.37'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.90'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.10'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.83'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.51'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +24.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.65'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.42%  '
$ws.Range('E32').Value = '  +5.76%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '680.07'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.82'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.21%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '66.67'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '42.74'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.87%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.432'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.158'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +6.38%  '
$ws.Range('D39').Value = '0.0₃0857'
$ws.Range('E39').Value = '  +3.14%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.46'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.43%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0498'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.19%  '
$ws.Range('E43').Value = '  -0.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.21'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.47%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.152'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +8.74%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.72'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.60%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.42'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.94%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.55'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +11.96%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.05'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.77%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.37'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.000273'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.21%  '
